# Applies the coin-price/volume refresh described by the commit diff.
# Columns D (Price) and E (Volume(1h)) store numeric-looking text, so we
# force a text NumberFormat before assigning the value - otherwise Excel
# auto-converts them to numbers/percentages and silently drops things like
# trailing zeros ("1.900" -> 1.9) or switches to scientific notation.
# Columns B (Coin) and C (Link) are plain text already, so a direct value
# assignment is enough.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '308.01'
# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '37.62'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '-0.18%'
# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '5.156'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '1.31%'
# Row 6
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '4.434'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '2.00%'
# Row 7
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.900'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '-0.16%'
# Row 8
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '8.244'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '0.44%'
# Row 9
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.991'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '4.22%'
# Row 10
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.9321'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '1.51%'
# Row 11
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.1067'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '-11.32%'
# Row 12
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.1923'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '0.09%'
# Row 13
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.08958'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '-2.38%'
# Row 14
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.03310'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '-2.59%'
# Row 15
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.09591'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '-1.02%'
# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.001387'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '0.44%'
# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.005918'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '3.91%'
# Row 18
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.609'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '1.51%'
# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.3388'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '0.38%'
# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.256'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '24.13%'
# Row 21
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '0.22%'
# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.2584'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '-0.03%'
# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.04397'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '0.93%'
# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.001232'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '1.80%'
# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.004559'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '6.99%'
# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.0001201'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '-7.61%'
# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.02195'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '3.72%'
# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.05042'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '1.02%'
# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.007456'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '-2.27%'
# Row 42
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '0.12%'
# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.008735'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '-10.91%'
# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.002111'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '2.65%'
# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.007983'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '-9.16%'
# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.00006536'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '-1.65%'
# Row 47
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '0.36%'
# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.002862'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '-5.75%'
# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.001002'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '-40.68%'
# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.00002104'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '0.36%'
# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0002003'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '0.36%'
